$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 1.73
$ws.Range("R2").Value = 2.1
$ws.Range("G5").Value = 2.1
$ws.Range("I5").Value = 4.2
$ws.Range("J5").Value = 3
$ws.Range("L5").Value = 5
$ws.Range("AD5").Value = 19
$ws.Range("O6").Value = 1.53
$ws.Range("P6").Value = 2.38
$ws.Range("Q6").Value = 2.7
$ws.Range("R6").Value = 1.44
$ws.Range("S6").Value = 4.5
$ws.Range("T6").Value = 1.21
$ws.Range("AR6").Value = 2.05
$ws.Range("AS6").Value = 1.8
$ws.Range("K9").Value = 2.2
$ws.Range("M9").Value = 1.08
$ws.Range("N9").Value = 8
$ws.Range("Q9").Value = 2.1
$ws.Range("R9").Value = 1.7
$ws.Range("S9").Value = 2.95
$ws.Range("T9").Value = 1.39
$ws.Range("U9").Value = 3.75
$ws.Range("V9").Value = 1.25
$ws.Range("Y9").Value = 2.63
$ws.Range("Z9").Value = 1.44
$ws.Range("AA9").Value = 5
$ws.Range("AE9").Value = 15
$ws.Range("AG9").Value = 8
$ws.Range("AH9").Value = 9
$ws.Range("AI9").Value = 29
$ws.Range("AL9").Value = 17
$ws.Range("AR9").Value = 1.56
$ws.Range("AS9").Value = 2.39
$ws.Range("G12").Value = 3.3
$ws.Range("I12").Value = 2.37
$ws.Range("J12").Value = 3.85
$ws.Range("K12").Value = 1.93
$ws.Range("L12").Value = 3
$ws.Range("N12").Value = 5.8
$ws.Range("U12").Value = 3.9
$ws.Range("W12").Value = 1.5
$ws.Range("X12").Value = 2.42
$ws.Range("AA12").Value = 8.5
$ws.Range("AB12").Value = 17
$ws.Range("AC12").Value = 11.25
$ws.Range("AD12").Value = 50
$ws.Range("AE12").Value = 32
$ws.Range("AG12").Value = 5.8
$ws.Range("AL12").Value = 6.5
$ws.Range("AM12").Value = 10.75
$ws.Range("AN12").Value = 9.25
$ws.Range("AO12").Value = 25
$ws.Range("AP12").Value = 22
$ws.Range("I16").Value = 5.5
$ws.Range("K16").Value = 2.05
$ws.Range("M16").Value = 1.08
$ws.Range("N16").Value = 8
$ws.Range("Q16").Value = 2.3
$ws.Range("R16").Value = 1.62
$ws.Range("U16").Value = 4.33
$ws.Range("V16").Value = 1.22
$ws.Range("W16").Value = 1.5
$ws.Range("X16").Value = 2.5
$ws.Range("Y16").Value = 2.2
$ws.Range("Z16").Value = 1.62
$ws.Range("AF16").Value = 34
$ws.Range("AG16").Value = 7.5
$ws.Range("AR16").Value = 1.78
$ws.Range("AS16").Value = 2.1
$ws.Range("H17").Value = 3.1
$ws.Range("L17").Value = 3.1
$ws.Range("N19").Value = 10
$ws.Range("L20").Value = 3.75
$ws.Range("W20").Value = 1.4
$ws.Range("X20").Value = 2.75
$ws.Range("Y20").Value = 1.73
$ws.Range("Z20").Value = 2
$ws.Range("AL20").Value = 10
$ws.Range("AQ20").Value = 34
$ws.Range("Y21").Value = 1.83
$ws.Range("Z21").Value = 1.83
$ws.Range("G22").Value = 2.3
$ws.Range("I22").Value = 2.88
$ws.Range("L22").Value = 3.4
$ws.Range("Y22").Value = 1.62
$ws.Range("AC22").Value = 9.5
$ws.Range("AM22").Value = 15
$ws.Range("AQ22").Value = 26
$ws.Range("G25").Value = 2.75
$ws.Range("I25").Value = 2.63
$ws.Range("J25").Value = 3.4
$ws.Range("K25").Value = 2.05
$ws.Range("L25").Value = 3.25
$ws.Range("N25").Value = 8.5
$ws.Range("Q25").Value = 2.08
$ws.Range("R25").Value = 1.73
$ws.Range("U25").Value = 3.75
$ws.Range("V25").Value = 1.25
$ws.Range("AA25").Value = 8.5
$ws.Range("AD25").Value = 26
$ws.Range("AG25").Value = 8.5
$ws.Range("AI25").Value = 15
$ws.Range("AM25").Value = 13
$ws.Range("AQ25").Value = 34
